# Populate the "rock-n-roll" vocabulary sheet with real content:
# fill in the metadata fields (rows 8-15) that previously contained
# placeholder/instruction text, flesh out the term rows 19-25 that only
# had a bare "rock-n-roll:" prefix, and append three new term rows
# (26-28) for Color, Green and a fresh blank template row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 'rock-n-roll'
$ws.Range("B9").Value = 'rock-n-roll'
$ws.Range("B10").Value = 'This Vocabulary describes some few terms of high level geologic concept hierarchy. This vocabulary is not complete and half randomly made up, but based on the USGS North American Geological Map Data Model. It is a test before the M4M workshop.'
$ws.Range("B11").Value = 'https://orcid.org/0000-0003-0682-8303'
$ws.Range("B12").Value = 'https://spdx.org/licenses/CC-BY-1.0.html'
$ws.Range("B13").Value = 'version 0.1'
$ws.Range("B14").Value = '2020-02-22T11:55:00+00:00'
$ws.Range("B15").Value = '2020-02-22T11:55:00+00:00'
$ws.Range("A19").Value = 'rock-n-roll:bestrock'
$ws.Range("B19").Value = 'bestrock'
$ws.Range("A20").Value = 'rock-n-roll:GeologicConcept'
$ws.Range("B20").Value = 'GeologicConcept'
$ws.Range("D20").Value = 'A subset of the Universe of all concepts, which includes only those concepts related to: "The study of the planet Earth--the materials of which it is made, the processes that act on these materials, the products formed, and the history of the planet and its life forms since its origin." [Jackson, 1997; p. 265] This conceptual model is focused primarily on geologic concepts that can be represented on geologic maps and diagrams. The term concept represents the notion of any mental phenomena that human beings use in their internal representation of the world. Webster’s dictionary [1996] uses the terms ‘idea’ and ‘object of thought’ to convey the meaning of ‘concept.’ GeologicConcepts identify the kinds of observable or inferred phenomena that earth scientists recognize.'
$ws.Range("A21").Value = 'rock-n-roll:GeologicProperty'
$ws.Range("B21").Value = 'GeologicProperty'
$ws.Range("D21").Value = 'An inherent feature used to characterize a GeologicConcept.'
$ws.Range("E21").Value = 'rock-n-roll:GeologicConcept'
$ws.Range("A22").Value = 'rock-n-roll:EarthMaterial'
$ws.Range("B22").Value = 'EarthMaterial'
$ws.Range("D22").Value = 'A naturally occuring substance in the Earth. EarthMaterialsare defined strictly by chemical and physical properties.'
$ws.Range("E22").Value = 'rock-n-roll:GeologicConcept'
$ws.Range("A23").Value = 'rock-n-roll:CompoundMaterial'
$ws.Range("B23").Value = 'CompoundMaterial'
$ws.Range("C23").Value = 'MineralMix'
$ws.Range("D23").Value = 'An EarthMaterial composed of other EarthMaterial instances, possibly including other CompoundMaterial instances. Includes consolidated and unconsolidated materials as well as mixtures of consolidated and unconsolidated materials.'
$ws.Range("E23").Value = 'rock-n-roll:EarthMaterial'
$ws.Range("A24").Value = 'rock-n-roll:Rock'
$ws.Range("B24").Value = 'Rock'
$ws.Range("C24").Value = 'Stone, ReallyHardThing'
$ws.Range("D24").Value = 'A consolidated aggregate of one or more EarthMaterials, or a body of undifferentiated mineral matter, or of solid organic material [adapted from Jackson, 1997]. Includes mineral aggregates such as granite, shale, marble; mineral matter, such as obsidian (Glass); and organic material, such a coal. Excludes unconsolidated materials.'
$ws.Range("E24").Value = 'rock-n-roll:CompundMaterial'
$ws.Range("A25").Value = 'rock-n-roll:Mineral'
$ws.Range("B25").Value = 'Mineral'
$ws.Range("D25").Value = 'A naturally occurring inorganic element or compound having a periodically repeating arrangement of atoms and a characteristic chemical composition or range of compositions, resulting in distinctive physical properties. Includes mercury as a general exception to the requirement of crystallinity. Also includes cryptocrystalline materials such as chalcedony and amorphous silica.'
$ws.Range("E25").Value = 'rock-n-roll:EarthMaterial'
$ws.Range("A26").Value = 'rock-n-roll:Color'
$ws.Range("B26").Value = 'Color'
$ws.Range("C26").Value = 'Colour'
$ws.Range("D26").Value = 'An attribute to describe the color of an instance of a concept such as EarthMaterial or GeologicUnit. A controlled vocabulary for color might be very useful.'
$ws.Range("E26").Value = 'rock-n-roll:GeologicProperty'
$ws.Range("A27").Value = 'rock-n-roll:Green'
$ws.Range("B27").Value = 'Green'
$ws.Range("C27").Value = 'Lime, DarkGreen, LightGreen'
$ws.Range("D27").Value = 'A color of an EarthMaterial'
$ws.Range("E27").Value = 'rock-n-roll:Color'
$ws.Range("A28").Value = 'rock-n-roll:'
